$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the columns that get reshuffled
# across rows 2..10, so the permutation can be applied safely regardless of
# the order in which we write the new values.
$cols = @("D","L","M","N","O","P","S")
$snapshot = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of new row -> source (old) row, derived from the diff.
$rowMap = @{
    2 = 5
    3 = 6
    4 = 2
    5 = 4
    6 = 10
    7 = 8
    8 = 9
    9 = 3
    10 = 7
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
